$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: IE users / IQ example
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "IE Users Have Lower IQ Than Users of Other Web Browsers [STUDY]"
$ws.Range("C5").Value = "technology, intelligence"
$ws.Range("D5").Value = "correlation"
$ws.Range("F5").Value = 2017
$ws.Range("G5").Value = "IQ"
$ws.Range("H5").Value = "web browser"
$ws.Range("L5").Value = 20181202

# Hyperlinks - K5 (source) first, then J5 (url), matching rId4/rId5 ordering
$ws.Hyperlinks.Add($ws.Range("K5"), "https://www.macmillanihe.com/companion/De-Vries-Critical-Statistics/")
$ws.Range("K5").Style = $ws.Range("J2").Style

$ws.Hyperlinks.Add($ws.Range("J5"), "https://mashable.com/2011/07/29/internet-explorer-iq/")
$ws.Range("J5").Style = $ws.Range("J2").Style

$ws.Range("A5").Select() | Out-Null

Write-Output "Applied edit"
